# Update the public EPEX Spot prices workbook with the new day column (08-jul)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# New header for column Y, matching style of existing headers (column X / row 1)
$ws.Range("X1").Copy()
$ws.Range("Y1").PasteSpecial(-4122)
$ws.Range("Y1").Value = "08-jul"

# New hourly price values for 08-jul
$values = @(61.85, 43.45, 39.17, 31.81, 31.91, 30.1, 38.52, 61.13, 60, 18.01, 4.31, 1.72, 1.72, 0, 0, 0, 0, 11.78, 45, 60.39, 86.87, 82.45999999999999, 101.59, 88.16)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 25).Value = $values[$i]
}
